$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to be treated as text so numeric-looking
# strings (e.g. "0.6276", "81.80") are not coerced into floating point
# numbers, matching the original inline-string cell content.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.375.79'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.841.16'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '239.13'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = '0.6276'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '0.07397'
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("D9").Value = '0.2891'
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").Value = '24.86'
$ws.Range("E10").Value = '  +1.88%  '
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '1.840.73'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '4.966'
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").Value = '0.6731'
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").Value = '0.00001021'
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").Value = '81.80'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '6.277'
$ws.Range("E17").Value = '  +1.61%  '
$ws.Range("D18").Value = '29.377.31'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '234.13'
$ws.Range("E19").Value = '  +2.61%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '7.297'
$ws.Range("E22").Value = '  -2.78%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '157.67'
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").Value = '8.493'
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").Value = '0.1342'
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("D27").Value = '17.30'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").Value = '0.07231'
$ws.Range("E28").Value = '  +11.71%  '
$ws.Range("D29").Value = '1.494'
$ws.Range("E29").Value = '  +5.01%  '
$ws.Range("D30").Value = '1.475'
$ws.Range("E30").Value = '  -0.59%  '
$ws.Range("D31").Value = '4.039'
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("D32").Value = '4.031'
$ws.Range("E32").Value = '  -1.38%  '
$ws.Range("D33").Value = '1.816'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("D35").Value = '0.6978'
$ws.Range("E35").Value = '  +0.71%  '
$ws.Range("D36").Value = '2.577'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").Value = '0.01829'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").Value = '2.805'
$ws.Range("E38").Value = '  -0.95%  '
$ws.Range("D39").Value = '1.232.62'
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("D40").Value = '6.765'
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("D41").Value = '0.9477'
$ws.Range("E41").Value = '  +1.70%  '
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").Value = '1.993.41'
$ws.Range("D44").Value = '101.11'
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("D45").Value = '65.28'
$ws.Range("E45").Value = '  -0.96%  '
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  +3.89%  '
$ws.Range("D47").Value = '1.697'
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").Value = '6.956'
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("D49").Value = '8.861'
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("D50").Value = '0.3888'
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("D51").Value = '0.1125'
$ws.Range("E51").Value = '  -2.40%  '

# Restore default (General) formatting so no stray cell style is left
# behind by the temporary text-number-format above.
$dataRange.ClearFormats()
